$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 119, shifting existing rows 119-127 down to 120-128.
$ws.Rows.Item(119).Insert()

# Populate the newly inserted row 119 with the new price record.
$ws.Cells.Item(119,1).Value = 8
$ws.Cells.Item(119,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(119,3).Value = "Coquimbo"
$ws.Cells.Item(119,4).Value = 44161
$ws.Cells.Item(119,5).Value = 4
$ws.Cells.Item(119,6).Value = "Fruta"
$ws.Cells.Item(119,7).Value = 100103
$ws.Cells.Item(119,8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(119,9).Value = 100103001
$ws.Cells.Item(119,10).Value = "Cereza"
$ws.Cells.Item(119,11).Value = "Rainier"
$ws.Cells.Item(119,12).Value = "Primera"
$ws.Cells.Item(119,13).Value = 300
$ws.Cells.Item(119,14).Value = 24500
$ws.Cells.Item(119,15).Value = 25000
$ws.Cells.Item(119,16).Value = 24750
$ws.Cells.Item(119,17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(119,18).Value = "Provincia de Curicó"
$ws.Cells.Item(119,19).Value = 2475
$ws.Cells.Item(119,20).Value = 10
